$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns.
# Numeric-looking price strings must be force-written as text (matching
# the source data which keeps thousands-separator formatted strings),
# so we briefly mark the cell as Text before assigning, then clear the
# format back off so no stray style index is left behind.

$ws.Range("D2").Value = "66.520.25"
$ws.Range("D3").Value = "2.585.44"
$ws.Range("E3").Value = "  -1.85%  "
$ws.Range("E4").Value = "  -0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "583.78"
$c.ClearFormats()
$ws.Range("E5").Value = "  -1.58%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "166.32"
$c.ClearFormats()
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -1.56%  "
$ws.Range("D9").Value = "2.583.88"
$ws.Range("E9").Value = "  -1.91%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.138"
$c.ClearFormats()
$ws.Range("E10").Value = "  -3.85%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("E13").Value = "  -1.62%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "26.76"
$c.ClearFormats()
$ws.Range("E14").Value = "  -3.87%  "
$ws.Range("D15").Value = "3.055.07"
$ws.Range("E16").Value = "  -2.49%  "
$ws.Range("D17").Value = "66.384.46"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "2.615.69"
$ws.Range("E18").Value = "  -0.16%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.44"
$c.ClearFormats()
$ws.Range("E19").Value = "  -6.03%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.73"
$c.ClearFormats()
$ws.Range("E20").Value = "  -4.10%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "352.25"
$c.ClearFormats()
$ws.Range("E22").Value = "  -2.95%  "
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  -3.92%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "68.77"
$c.ClearFormats()
$ws.Range("E26").Value = "  -2.39%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.95"
$c.ClearFormats()
$ws.Range("E27").Value = "  -8.67%  "
$ws.Range("D28").Value = "2.715.81"
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("D29").Value = "0.0₃0987"
$ws.Range("E29").Value = "  -2.76%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "533.31"
$c.ClearFormats()
$ws.Range("E30").Value = "  -4.21%  "
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("E33").Value = "  -3.06%  "
$ws.Range("E34").Value = "  -3.04%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -3.61%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "156.93"
$c.ClearFormats()
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("E38").Value = "  -2.37%  "
$ws.Range("E39").Value = "  -1.80%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "18.26"
$c.ClearFormats()
$ws.Range("E40").Value = "  +1.78%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.78"
$c.ClearFormats()
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("E44").Value = "  -2.60%  "
$ws.Range("E45").Value = "  -4.85%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "149.40"
$c.ClearFormats()
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("E47").Value = "  -3.63%  "
$ws.Range("E48").Value = "  -2.58%  "
$ws.Range("E49").Value = "  -1.46%  "
$ws.Range("E50").Value = "  -1.66%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.596"
$c.ClearFormats()
$ws.Range("E51").Value = "  -1.44%  "
